$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 35365
$ws.Range("B6").Value = "kjlk"
$ws.Range("C6").Value = 8422458254
$ws.Range("D6").Value = 44387
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E12").Select()
